$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 75.92737354118466
$ws.Range("R2").Value = 683.346361870662
$ws.Range("S2").Value = 0.005858992042618862
$ws.Range("T2").Value = 0.005858992042618862
$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 126.6369783078867
$ws.Range("R3").Value = 1139.73280477098
$ws.Range("S3").Value = 0.009772036271012948
$ws.Range("T3").Value = 0.009772036271012948
$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 13.71176749813267
$ws.Range("R4").Value = 123.405907483194
$ws.Range("S4").Value = 0.00105807869961711
$ws.Range("T4").Value = 0.00105807869961711
$ws.Range("G5").Value = 29.223446
$ws.Range("H5").Value = 87.670338
$ws.Range("I5").Value = 0.0169041244192178
$ws.Range("J5").Value = 0.0169041244192178
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 2.786436093802666
$ws.Range("R5").Value = 25.077924844224
$ws.Range("S5").Value = 0.0002150174059688814
$ws.Range("T5").Value = 0.0002150174059688814
$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 4254.090350756108
$ws.Range("R6").Value = 38286.81315680497
$ws.Range("S6").Value = 0.328270034260332
$ws.Range("T6").Value = 0.328270034260332
$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 7095.26910181181
$ws.Range("S7").Value = 0.5475116979412674
$ws.Range("T7").Value = 0.5475116979412673
$ws.Range("I8").Value = 0.9471112884046843
$ws.Range("J8").Value = 0.9471112884046842
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 768.2485918464878
$ws.Range("R8").Value = 6914.237326618389
$ws.Range("S8").Value = 0.05928247187346984
$ws.Range("T8").Value = 0.05928247187346983
$ws.Range("I9").Value = 0.9471112884046843
$ws.Range("J9").Value = 0.9471112884046842
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 156.1195962246044
$ws.Range("R9").Value = 1405.07636602144
$ws.Range("S9").Value = 0.01204708432961496
$ws.Range("T9").Value = 0.01204708432961495
$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 97.1509646370689
$ws.Range("R10").Value = 874.3586817336201
$ws.Range("S10").Value = 0.00749672617652951
$ws.Range("T10").Value = 0.007496726176529508
$ws.Range("G11").Value = 37.39212666666667
$ws.Range("H11").Value = 112.17638
$ws.Range("I11").Value = 0.02162924801792661
$ws.Range("J11").Value = 0.0216292480179266
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 162.0351663377556
$ws.Range("R11").Value = 1458.3164970398
$ws.Range("S11").Value = 0.01250356368092172
$ws.Range("T11").Value = 0.01250356368092172
$ws.Range("G12").Value = 37.39212666666667
$ws.Range("H12").Value = 112.17638
$ws.Range("I12").Value = 0.02162924801792661
$ws.Range("J12").Value = 0.0216292480179266
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 17.54454786454889
$ws.Range("R12").Value = 157.90093078094
$ws.Range("S12").Value = 0.001353838036738889
$ws.Range("T12").Value = 0.001353838036738889
$ws.Range("G13").Value = 37.39212666666667
$ws.Range("H13").Value = 112.17638
$ws.Range("I13").Value = 0.02162924801792661
$ws.Range("J13").Value = 0.0216292480179266
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 3.565314349582222
$ws.Range("R13").Value = 32.08782914624
$ws.Range("S13").Value = 0.0002751201237364856
$ws.Range("T13").Value = 0.0002751201237364855
$ws.Range("G14").Value = 24.817167
$ws.Range("H14").Value = 74.45150100000001
$ws.Range("I14").Value = 0.01435533915817136
$ws.Range("J14").Value = 0.01435533915817136
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 64.479127788111
$ws.Range("R14").Value = 580.3121500929991
$ws.Range("S14").Value = 0.004975579675762517
$ws.Range("T14").Value = 0.004975579675762516
$ws.Range("G15").Value = 24.817167
$ws.Range("H15").Value = 74.45150100000001
$ws.Range("I15").Value = 0.01435533915817136
$ws.Range("J15").Value = 0.01435533915817136
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 107.54279420169
$ws.Range("R15").Value = 967.8851478152101
$ws.Range("S15").Value = 0.0082986194053838
$ws.Range("T15").Value = 0.0082986194053838
$ws.Range("G16").Value = 24.817167
$ws.Range("H16").Value = 74.45150100000001
$ws.Range("I16").Value = 0.01435533915817136
$ws.Range("J16").Value = 0.01435533915817136
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 11.644322297457
$ws.Range("R16").Value = 104.798900677113
$ws.Range("S16").Value = 0.0008985427586993218
$ws.Range("T16").Value = 0.0008985427586993217
$ws.Range("G17").Value = 24.817167
$ws.Range("H17").Value = 74.45150100000001
$ws.Range("I17").Value = 0.01435533915817136
$ws.Range("J17").Value = 0.01435533915817136
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 2.366300328672
$ws.Range("R17").Value = 21.296702958048
$ws.Range("S17").Value = 0.0002751201237364856
$ws.Range("T17").Value = 0.0002751201237364855

Write-Output "done"
